$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# Status column (I) for all incident rows moved from "Pendente" to "Resolvido"
$ws.Range("I2:I22").Value = "Resolvido"

# Reflect the final on-screen selection recorded for this sheet
$ws.Activate()
$ws.Range("M17:M18").Select() | Out-Null
